$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column before column A. This pushes every existing
# column (id, title, display_name, ... domain) one slot to the right,
# i.e. the old A:AP range becomes B:AQ.
$ws.Range("A1").EntireColumn.Insert()

# Header for the new column, styled the same bold/centered way as the
# rest of row 1 (reuses the existing header style, no new style needed).
$ws.Range("A1").Value = "UA_authored_year"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108

# Data value for row 2. Force text storage (so it matches "2023" as a
# string, not a number) by switching the cell to a text number format
# before writing, then clearing the format override afterwards so the
# cell ends up with no explicit style, same as its neighbours.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2023"
$ws.Range("A2").ClearFormats()
